# Rebuild the word-frequency table (Sheet1) with the new word list and counts.
# ClearContents() wipes cell values but keeps the existing cell styles (s="1")
# intact, and writing the words back in row order gives the shared-string table
# a fresh, deterministic rebuild (matching a from-scratch recompute/save).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.ClearContents()

$ws.Range("B1").Value = "단어"
$ws.Range("C1").Value = "빈도수"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "버스"
$ws.Range("C2").Value = 1492
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "퇴근"
$ws.Range("C3").Value = 763
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "선"
$ws.Range("C4").Value = 752
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "노선"
$ws.Range("C5").Value = 651
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "동"
$ws.Range("C6").Value = 570
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "분"
$ws.Range("C7").Value = 464
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "출근"
$ws.Range("C8").Value = 434
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "없"
$ws.Range("C9").Value = 422
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "시간"
$ws.Range("C10").Value = 396
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "사람"
$ws.Range("C11").Value = 366
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "영통"
$ws.Range("C12").Value = 363
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "이용"
$ws.Range("C13").Value = 361
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "차"
$ws.Range("C14").Value = 337
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "많"
$ws.Range("C15").Value = 333
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "증차"
$ws.Range("C16").Value = 289
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "운행"
$ws.Range("C17").Value = 281
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "잠실"
$ws.Range("C18").Value = 272
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "경우"
$ws.Range("C19").Value = 270
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "인원"
$ws.Range("C20").Value = 260
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "생각"
$ws.Range("C21").Value = 251
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "차량"
$ws.Range("C22").Value = 249
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "시"
$ws.Range("C23").Value = 237
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "도착"
$ws.Range("C24").Value = 231
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "출발"
$ws.Range("C25").Value = 224
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "자리"
$ws.Range("C26").Value = 222
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "요청"
$ws.Range("C27").Value = 209
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "건의"
$ws.Range("C28").Value = 205
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "부탁"
$ws.Range("C29").Value = 201
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "회사"
$ws.Range("C30").Value = 181
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "분당"
$ws.Range("C31").Value = 181
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "감사"
$ws.Range("C32").Value = 175
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "기사님"
$ws.Range("C33").Value = 173
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "검토"
$ws.Range("C34").Value = 169
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "정도"
$ws.Range("C35").Value = 161
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "금요일"
$ws.Range("C36").Value = 152
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "좋"
$ws.Range("C37").Value = 151
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "통근"
$ws.Range("C38").Value = 150
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "문제"
$ws.Range("C39").Value = 149
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "불편"
$ws.Range("C40").Value = 140
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "아침"
$ws.Range("C41").Value = 138
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "출퇴근"
$ws.Range("C42").Value = 138
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "글"
$ws.Range("C43").Value = 137
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "안녕"
$ws.Range("C44").Value = 135
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "배차"
$ws.Range("C45").Value = 134
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = "관련"
$ws.Range("C46").Value = 130
$ws.Range("A47").Value = 45
$ws.Range("B47").Value = "답변"
$ws.Range("C47").Value = 128
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = "추가"
$ws.Range("C48").Value = 124
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = "정류장"
$ws.Range("C49").Value = 120
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = "고속도로"
$ws.Range("C50").Value = 119
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = "상황"
$ws.Range("C51").Value = 116
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = "변경"
$ws.Range("C52").Value = 116
$ws.Range("A53").Value = 51
$ws.Range("B53").Value = "양재"
$ws.Range("C53").Value = 116
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = "일"
$ws.Range("C54").Value = 112
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = "등"
$ws.Range("C55").Value = 111
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = "마을"
$ws.Range("C56").Value = 111
$ws.Range("A57").Value = 55
$ws.Range("B57").Value = "역"
$ws.Range("C57").Value = 109
$ws.Range("A58").Value = 56
$ws.Range("B58").Value = "사항"
$ws.Range("C58").Value = 107
$ws.Range("A59").Value = 57
$ws.Range("B59").Value = "필요"
$ws.Range("C59").Value = 105
$ws.Range("A60").Value = 58
$ws.Range("B60").Value = "앞"
$ws.Range("C60").Value = 105
$ws.Range("A61").Value = 59
$ws.Range("B61").Value = "직원"
$ws.Range("C61").Value = 102
$ws.Range("A62").Value = 60
$ws.Range("B62").Value = "안"
$ws.Range("C62").Value = 99
$ws.Range("A63").Value = 61
$ws.Range("B63").Value = "쪽"
$ws.Range("C63").Value = 96
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = "현재"
$ws.Range("C64").Value = 96
$ws.Range("A65").Value = 63
$ws.Range("B65").Value = "행"
$ws.Range("C65").Value = 95
$ws.Range("A66").Value = 64
$ws.Range("B66").Value = "전"
$ws.Range("C66").Value = 94
$ws.Range("A67").Value = 65
$ws.Range("B67").Value = "안양"
$ws.Range("C67").Value = 92
$ws.Range("A68").Value = 66
$ws.Range("B68").Value = "개선"
$ws.Range("C68").Value = 91
$ws.Range("A69").Value = 67
$ws.Range("B69").Value = "중"
$ws.Range("C69").Value = 91
$ws.Range("A70").Value = 68
$ws.Range("B70").Value = "탑승"
$ws.Range("C70").Value = 91
$ws.Range("A71").Value = 69
$ws.Range("B71").Value = "번"
$ws.Range("C71").Value = 90
$ws.Range("A72").Value = 70
$ws.Range("B72").Value = "수원"
$ws.Range("C72").Value = 90
$ws.Range("A73").Value = 71
$ws.Range("B73").Value = "성남"
$ws.Range("C73").Value = 88
$ws.Range("A74").Value = 72
$ws.Range("B74").Value = "이상"
$ws.Range("C74").Value = 87
$ws.Range("A75").Value = 73
$ws.Range("B75").Value = "길"
$ws.Range("C75").Value = 87
$ws.Range("A76").Value = 74
$ws.Range("B76").Value = "때문"
$ws.Range("C76").Value = 87
$ws.Range("A77").Value = 75
$ws.Range("B77").Value = "운영"
$ws.Range("C77").Value = 86
$ws.Range("A78").Value = 76
$ws.Range("B78").Value = "중부"
$ws.Range("C78").Value = 86
$ws.Range("A79").Value = 77
$ws.Range("B79").Value = "조정"
$ws.Range("C79").Value = 85
$ws.Range("A80").Value = 78
$ws.Range("B80").Value = "금"
$ws.Range("C80").Value = 85
$ws.Range("A81").Value = 79
$ws.Range("B81").Value = "수지"
$ws.Range("C81").Value = 83
$ws.Range("A82").Value = 80
$ws.Range("B82").Value = "경유"
$ws.Range("C82").Value = 80
$ws.Range("A83").Value = 81
$ws.Range("B83").Value = "구리"
$ws.Range("C83").Value = 79
$ws.Range("A84").Value = 82
$ws.Range("B84").Value = "초과"
$ws.Range("C84").Value = 78
$ws.Range("A85").Value = 83
$ws.Range("B85").Value = "신설"
$ws.Range("C85").Value = 77
$ws.Range("A86").Value = 84
$ws.Range("B86").Value = "발생"
$ws.Range("C86").Value = 76
$ws.Range("A87").Value = 85
$ws.Range("B87").Value = "오늘"
$ws.Range("C87").Value = 76
$ws.Range("A88").Value = 86
$ws.Range("B88").Value = "대"
$ws.Range("C88").Value = 74
$ws.Range("A89").Value = 87
$ws.Range("B89").Value = "운전"
$ws.Range("C89").Value = 74
$ws.Range("A90").Value = 88
$ws.Range("B90").Value = "그렇"
$ws.Range("C90").Value = 73
$ws.Range("A91").Value = 89
$ws.Range("B91").Value = "문의"
$ws.Range("C91").Value = 72
$ws.Range("A92").Value = 90
$ws.Range("B92").Value = "천호"
$ws.Range("C92").Value = 72
$ws.Range("A93").Value = 91
$ws.Range("B93").Value = "후"
$ws.Range("C93").Value = 71
$ws.Range("A94").Value = 92
$ws.Range("B94").Value = "담당자"
$ws.Range("C94").Value = 71
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = "좌석"
$ws.Range("C95").Value = 69
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = "최근"
$ws.Range("C96").Value = 68
$ws.Range("A97").Value = 95
$ws.Range("B97").Value = "광주"
$ws.Range("C97").Value = 68
$ws.Range("A98").Value = 96
$ws.Range("B98").Value = "사고"
$ws.Range("C98").Value = 68
$ws.Range("A99").Value = 97
$ws.Range("B99").Value = "부분"
$ws.Range("C99").Value = 67
$ws.Range("A100").Value = 98
$ws.Range("B100").Value = "부족"
$ws.Range("C100").Value = 67
$ws.Range("A101").Value = 99
$ws.Range("B101").Value = "어떻"
$ws.Range("C101").Value = 67
